# Update regression test data for preprod R33 run.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Switch the environment/account under test from the "preproducciongestion"
# host to the "i-preproducciongestion" host, and roll the account number.
$ws.Range("A2").Value = "i-preproducciongestion.segurossura.com.ar"
$ws.Range("B2").Value = "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"
$ws.Range("E2").Value = 5069929970

# Leave the sheet scrolled/selected the way it was when the data was last
# captured (selection on S2, scrolled so column L is at the left edge).
$ws.Range("S2").Select()
$excel.ActiveWindow.ScrollColumn = 12
